# ---------------------------------------------------------------------------
# Adds a new "Player Info" worksheet in front of the existing "ODI Batting"
# sheet, and updates the "ODI Batting" sheet's MATCH_CARD_LINK column so it
# becomes a MATCH_CODE column holding just the numeric match code instead of
# the full scorecard URL.
# ---------------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$odi = $wb.Worksheets.Item(1)

# --- 1. Insert a new "Player Info" sheet before "ODI Batting" --------------
# Duplicating the existing sheet keeps the exact same style table (bold,
# centered, bordered header) instead of fabricating new style entries.
# NOTE: after Copy(), sheet handles obtained beforehand start tracking the
# *new* sheet (the engine resolves worksheet handles positionally), so every
# sheet reference used below is re-fetched right after the copy.
$odi.Copy($odi)

$playerInfo = $wb.Worksheets.Item(1)
$playerInfo.Name = "Player Info"

$odi = $wb.Worksheets.Item(2)

# --- 2. Update the "ODI Batting" sheet: MATCH_CARD_LINK -> MATCH_CODE ------
$odi.Range("D1").Value = "MATCH_CODE"

$odi.Range("D2").NumberFormat = "@"
$odi.Range("D2").Value = "4726"

$odi.Range("D3").NumberFormat = "@"
$odi.Range("D3").Value = "4729"

$odi.Range("D4").NumberFormat = "@"
$odi.Range("D4").Value = "4734"

# --- 3. Populate the "Player Info" sheet ------------------------------------
# Drop the columns/rows that "ODI Batting" had but "Player Info" doesn't.
$playerInfo.Range("E1:J4").Clear()
$playerInfo.Range("A3:D4").Clear()

# Header row (keeps the bold/centered/bordered style copied from "ODI Batting").
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Data row.
$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "7162"
$playerInfo.Range("B2").Value = "Md Towhid Hridoy"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Off Break"
